$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price (D) and volume-change (E) values
$ws.Range('D2').Value = '26.771.46'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').Value = '1.536.25'
$ws.Range('E3').Value = '  -1.82%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '205.27'
$ws.Range('E5').Value = '  -0.46%  '
$ws.Range('E6').Value = '  -1.01%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -0.95%  '
$ws.Range('E9').Value = '  -2.90%  '
$ws.Range('E10').Value = '  -0.56%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0854'
$ws.Range('E11').Value = '  -0.89%  '
$ws.Range('D12').Value = '1.753.90'
$ws.Range('E12').Value = '  -1.89%  '
$ws.Range('D13').Value = '1.531.28'
$ws.Range('E13').Value = '  -2.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.66'
$ws.Range('E14').Value = '  -1.55%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.507'
$ws.Range('E15').Value = '  -1.40%  '
$ws.Range('D16').Value = '26.762.63'
$ws.Range('E16').Value = '  -0.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '60.89'
$ws.Range('E17').Value = '  -1.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '212.70'
$ws.Range('E18').Value = '  -0.60%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.22'
$ws.Range('E19').Value = '  -1.57%  '
$ws.Range('E20').Value = '  +0.63%  '
$ws.Range('E22').Value = '  -1.92%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.09'
$ws.Range('E23').Value = '  -2.42%  '
$ws.Range('E24').Value = '  -3.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.24'
$ws.Range('E25').Value = '  -0.74%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.57'
$ws.Range('E26').Value = '  -2.34%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '14.76'
$ws.Range('E27').Value = '  -0.91%  '
$ws.Range('E28').Value = '  -0.11%  '
$ws.Range('E29').Value = '  -1.06%  '
$ws.Range('E30').Value = '  -1.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0455'
$ws.Range('E31').Value = '  -1.55%  '
$ws.Range('E32').Value = '  +2.35%  '
$ws.Range('D33').Value = '1.362.53'
$ws.Range('E33').Value = '  -1.63%  '
$ws.Range('E34').Value = '  -0.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.50'
$ws.Range('E35').Value = '  -2.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.956'
$ws.Range('E36').Value = '  +3.41%  '
$ws.Range('E38').Value = '  +0.58%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.520'
$ws.Range('E39').Value = '  -0.49%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.73'
$ws.Range('E40').Value = '  +7.44%  '
$ws.Range('E41').Value = '  -1.75%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.995'
$ws.Range('E42').Value = '  +0.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.19'
$ws.Range('E43').Value = '  +0.49%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '62.66'
$ws.Range('E44').Value = '  -0.85%  '
$ws.Range('E45').Value = '  -3.36%  '
$ws.Range('D46').Value = '1.669.10'
$ws.Range('E46').Value = '  -1.87%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '85.08'
$ws.Range('E47').Value = '  -0.33%  '
$ws.Range('E48').Value = '  +2.54%  '
$ws.Range('D49').Value = '0.0₇0980'
$ws.Range('E49').Value = '  -0.54%  '
$ws.Range('E50').Value = '  -0.89%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.999'
$ws.Range('E51').Value = '  -0.09%  '
